# Grainger.xlsx: rename sheet and update the saved selection/scroll state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomeDepot_URL")

# Sheet tab name: "HomeDepot_URL" -> "URL"
$ws.Name = "URL"

# Bring the sheet to the front and move the selection to D31.
# Selecting a cell in the (currently scrolled-to-column-J) view also resets
# the window's top-left cell back to the default (A1), matching the diff's
# removal of topLeftCell="J1".
$ws.Activate()
$ws.Range("D31").Select()
